$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update parameter values (B6:B8)
$ws.Range("B6").Value = 0.111
$ws.Range("B7").Value = 3.0550000000000002
$ws.Range("B8").Value = 2.99

# Update the active selection to B9 (matches diff's selection change)
$ws.Range("B9").Select()
